# Generate Report for Handback
#
# The handback pipeline produced a new result for file
# a786355d-db0b-41b5-8b60-ec3996eb278f: the handback transform failed
# (file name mismatch). This promotes that file's row to the top of each
# sheet (row 2), pushes 3f5144df-...md down to row 3, records the new
# status/date/error-detail, and leaves fddb60c5-...md (row 4) untouched.

$wb = $excel.ActiveWorkbook

function Set-Display {
    param($ws, [string]$cellAddr, [string]$display)
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq $cellAddr) {
            $hl.TextToDisplay = $display
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "a786355d-db0b-41b5-8b60-ec3996eb278f.md"
$wsOverview.Range("B2").Value = "Handback transform failed"
$wsOverview.Range("C2").Value = "Handback transform failed"
$wsOverview.Range("D2").Value = "2016-48-13 10:48:34"

$wsOverview.Range("A3").Value = "3f5144df-f6e2-4e3f-810b-022f6aa86f63.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-48-13 10:48:00"

Set-Display $wsOverview '$A$2' "a786355d-db0b-41b5-8b60-ec3996eb278f.md"
Set-Display $wsOverview '$A$3' "3f5144df-f6e2-4e3f-810b-022f6aa86f63.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "a786355d-db0b-41b5-8b60-ec3996eb278f.md"
$wsZh.Range("C2").Value = "Handback transform failed"
$wsZh.Range("D2").Value = "a786355d-db0b-41b5-8b60-ec3996eb278f.631357a838623148c02f593d860607143f7293c4.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-13 10:48:31"
$wsZh.Range("K2").Value = "Handback file name: iqi3cb45.cyr is different with handoff file name: a786355d-db0b-41b5-8b60-ec3996eb278f.631357a838623148c02f593d860607143f7293c4.zh-cn."

$wsZh.Range("A3").Value = "3f5144df-f6e2-4e3f-810b-022f6aa86f63.md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "3f5144df-f6e2-4e3f-810b-022f6aa86f63.0ad28e69cc55de91b9f7b12ea8609d2b235fb5e6.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-13 10:46:18"

Set-Display $wsZh '$A$2' "a786355d-db0b-41b5-8b60-ec3996eb278f.md"
Set-Display $wsZh '$D$2' "a786355d-db0b-41b5-8b60-ec3996eb278f.631357a838623148c02f593d860607143f7293c4.zh-cn.xlf"
Set-Display $wsZh '$A$3' "3f5144df-f6e2-4e3f-810b-022f6aa86f63.md"
Set-Display $wsZh '$D$3' "3f5144df-f6e2-4e3f-810b-022f6aa86f63.0ad28e69cc55de91b9f7b12ea8609d2b235fb5e6.zh-cn.xlf"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "a786355d-db0b-41b5-8b60-ec3996eb278f.md"
$wsDe.Range("C2").Value = "Handback transform failed"
$wsDe.Range("D2").Value = "a786355d-db0b-41b5-8b60-ec3996eb278f.631357a838623148c02f593d860607143f7293c4.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-13 10:48:34"
$wsDe.Range("K2").Value = "Handback file name: iqi3cb45.cyr is different with handoff file name: a786355d-db0b-41b5-8b60-ec3996eb278f.631357a838623148c02f593d860607143f7293c4.de-de."

$wsDe.Range("A3").Value = "3f5144df-f6e2-4e3f-810b-022f6aa86f63.md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "3f5144df-f6e2-4e3f-810b-022f6aa86f63.0ad28e69cc55de91b9f7b12ea8609d2b235fb5e6.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-13 10:48:00"

Set-Display $wsDe '$A$2' "a786355d-db0b-41b5-8b60-ec3996eb278f.md"
Set-Display $wsDe '$D$2' "a786355d-db0b-41b5-8b60-ec3996eb278f.631357a838623148c02f593d860607143f7293c4.de-de.xlf"
Set-Display $wsDe '$A$3' "3f5144df-f6e2-4e3f-810b-022f6aa86f63.md"
Set-Display $wsDe '$D$3' "3f5144df-f6e2-4e3f-810b-022f6aa86f63.0ad28e69cc55de91b9f7b12ea8609d2b235fb5e6.de-de.xlf"
